$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and restore WrappedEther/Polkadot row order)
$ws.Range('D2').Value = '25.885.79'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '1.639.16'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('D4').Value = '''1.003'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''215.62'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('D6').Value = '''0.5080'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '''0.2594'
$ws.Range('E8').Value = '  +1.06%  '
$ws.Range('D9').Value = '''0.06455'
$ws.Range('E9').Value = '  +1.60%  '
$ws.Range('D10').Value = '''20.24'
$ws.Range('E10').Value = '  +5.04%  '
$ws.Range('D11').Value = '''0.07824'
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = '''4.266'
$ws.Range('E12').Value = '  +0.58%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.645.38'
$ws.Range('E13').Value = '  +1.45%  '
$ws.Range('D14').Value = '1.865.66'
$ws.Range('E14').Value = '  +1.23%  '
$ws.Range('D15').Value = '''0.5665'
$ws.Range('E15').Value = '  +2.39%  '
$ws.Range('D16').Value = '0.0₅7692'
$ws.Range('E16').Value = '  +2.41%  '
$ws.Range('D17').Value = '''63.38'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').Value = '25.898.32'
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').Value = '''194.39'
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('D21').Value = '''4.393'
$ws.Range('E21').Value = '  +1.39%  '
$ws.Range('D22').Value = '''9.981'
$ws.Range('E22').Value = '  +2.41%  '
$ws.Range('D23').Value = '''6.235'
$ws.Range('E23').Value = '  +4.55%  '
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('D25').Value = '''1.756'
$ws.Range('E25').Value = '  -4.09%  '
$ws.Range('D26').Value = '''138.64'
$ws.Range('E26').Value = '  -1.42%  '
$ws.Range('E27').Value = '  -2.56%  '
$ws.Range('D28').Value = '''6.854'
$ws.Range('E28').Value = '  +2.06%  '
$ws.Range('D29').Value = '''15.55'
$ws.Range('E29').Value = '  +0.96%  '
$ws.Range('D30').Value = '''1.243'
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('D31').Value = '''0.04972'
$ws.Range('E31').Value = '  +2.43%  '
$ws.Range('D32').Value = '''3.320'
$ws.Range('E32').Value = '  +0.59%  '
$ws.Range('D33').Value = '''3.257'
$ws.Range('E33').Value = '  +2.67%  '
$ws.Range('D34').Value = '''1.576'
$ws.Range('E34').Value = '  +1.74%  '
$ws.Range('E35').Value = '  +1.15%  '
$ws.Range('D36').Value = '''0.9093'
$ws.Range('E36').Value = '  +1.90%  '
$ws.Range('D37').Value = '''2.575'
$ws.Range('E37').Value = '  +1.71%  '
$ws.Range('D38').Value = '''0.5527'
$ws.Range('E38').Value = '  +1.05%  '
$ws.Range('D39').Value = '1.127.11'
$ws.Range('E39').Value = '  +0.16%  '
$ws.Range('D40').Value = '''0.01573'
$ws.Range('E40').Value = '  +0.93%  '
$ws.Range('D41').Value = '''0.9999'
$ws.Range('E41').Value = '  -0.65%  '
$ws.Range('D42').Value = '''5.504'
$ws.Range('E42').Value = '  -1.10%  '
$ws.Range('D43').Value = '''99.53'
$ws.Range('E43').Value = '  +2.54%  '
$ws.Range('D44').Value = '''0.8014'
$ws.Range('E44').Value = '  +1.34%  '
$ws.Range('D45').Value = '0.0₈111'
$ws.Range('E45').Value = '  -1.44%  '
$ws.Range('D46').Value = '''55.64'
$ws.Range('E46').Value = '  +1.94%  '
$ws.Range('D47').Value = '''0.4239'
$ws.Range('E47').Value = '  -4.09%  '
$ws.Range('D48').Value = '''0.05044'
$ws.Range('E48').Value = '  -0.40%  '
$ws.Range('D49').Value = '''7.677'
$ws.Range('E49').Value = '  +2.06%  '
$ws.Range('D50').Value = '''0.9994'
$ws.Range('E50').Value = '  -0.12%  '
$ws.Range('D51').Value = '''1.004'
$ws.Range('E51').Value = '  +0.20%  '
